$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "btech"

# Header row (row 1)
$ws.Range('A1').Value = 'University_RollNumber'
$ws.Range('B1').Value = 'First_Name'
$ws.Range('C1').Value = 'Last_Name'
$ws.Range('D1').Value = 'Gender'
$ws.Range('E1').Value = 'Nationality'
$ws.Range('F1').Value = 'DOB'
$ws.Range('G1').Value = 'Phone_Number'
$ws.Range('H1').Value = 'Email_ID'
$ws.Range('I1').Value = 'ADHAR_Number'
$ws.Range('J1').Value = 'Address'
$ws.Range('K1').Value = 'District'
$ws.Range('L1').Value = 'State'
$ws.Range('M1').Value = 'Country'
$ws.Range('N1').Value = 'Pin_Code'
$ws.Range('O1').Value = 'Category'
$ws.Range('P1').Value = 'Sub_Category'
$ws.Range('Q1').Value = '_10th_CGPA'
$ws.Range('R1').Value = '_10th_Board'
$ws.Range('S1').Value = '_10th_YOP'
$ws.Range('T1').Value = '_12th_Percentage'
$ws.Range('U1').Value = '_12th_Board'
$ws.Range('V1').Value = '_12th_YOP'
$ws.Range('W1').Value = 'Diploma_Percentage'
$ws.Range('X1').Value = 'Diploma_Board'
$ws.Range('Y1').Value = 'Diploma_YOP'
$ws.Range('Z1').Value = 'Course_RegularORIntegrated'
$ws.Range('AA1').Value = 'Branch'
$ws.Range('AB1').Value = 'College_Name'
$ws.Range('AC1').Value = 'Course_CGPA'
$ws.Range('AD1').Value = 'Number_Of_Backlogs'
$ws.Range('AE1').Value = 'Entrance_Exam'
$ws.Range('AF1').Value = 'CET_Rank'
$ws.Range('AG1').Value = 'Course_YOP'
$ws.Range('AH1').Value = 'Certificate_Course'
$ws.Range('AI1').Value = 'Certificate_IssuedBy'
$ws.Range('AJ1').Value = 'CertificateUpload'
$ws.Range('AK1').Value = 'Program_name'
$ws.Range('AL1').Value = 'Program_code'
$ws.Range('AM1').Value = 'InternUpload'
$ws.Range('AN1').Value = 'CertificatePlatform'
$ws.Range('AO1').Value = 'StudyingYear'
$ws.Range('AP1').Value = 'Year'
$ws.Range('AQ1').Value = 'Name_of_the_Teacher'
$ws.Range('AR1').Value = 'Contact_Details'
$ws.Range('AS1').Value = 'Program_graduated_from'
$ws.Range('AT1').Value = 'Name_of_company'
$ws.Range('AU1').Value = 'Name_of_employer_with_contact_details'
$ws.Range('AV1').Value = 'Pay_Package_at_appointment'
$ws.Range('AW1').Value = 'NameOfTeacher'
$ws.Range('AX1').Value = 'Name_Of_Students'
$ws.Range('AY1').Value = 'Name_Of_Institution_joined'
$ws.Range('AZ1').Value = 'Name_Of_Programme_Admitted_To'
$ws.Range('BA1').Value = 'Upload'
$ws.Range('BB1').Value = 'list_of_students_undertaking'
$ws.Range('BC1').Value = 'Program_Graduated'

# Data row (row 2)
$ws.Range('A2').Value = "'421206421015"
$ws.Range('B2').Value = "'Jaya"
$ws.Range('C2').Value = "'Madhuri"
$ws.Range('D2').Value = "'female"
$ws.Range('E2').Value = "'Indian"
$ws.Range('F2').Value = "'1-23-2001"
$ws.Range('G2').Value = "'891987"
$ws.Range('H2').Value = "'421206421015@andhrauniversity.edu.in"
$ws.Range('I2').Value = "'9630798"
$ws.Range('J2').Value = "'14/205-A"
$ws.Range('K2').Value = "'Krishna"
$ws.Range('L2').Value = "'AP"
$ws.Range('M2').Value = "'India"
$ws.Range('N2').Value = "'521301"
$ws.Range('O2').Value = "'BC-B"
$ws.Range('P2').Value = "'padmasali"
$ws.Range('Q2').Value = "'9"
$ws.Range('R2').Value = "'ssc"
$ws.Range('S2').Value = 2016
$ws.Range('T2').Value = "'98"
$ws.Range('U2').Value = "'IPE"
$ws.Range('V2').Value = 2018
$ws.Range('W2').Value = "'"
$ws.Range('X2').Value = "'"
$ws.Range('Y2').Value = "'"
$ws.Range('Z2').Value = "'Regular"
$ws.Range('AA2').Value = "'CSE"
$ws.Range('AB2').Value = "'AUCE"
$ws.Range('AC2').Value = "'9"
$ws.Range('AD2').Value = 0
$ws.Range('AE2').Value = "'AUCET"
$ws.Range('AF2').Value = 1872
$ws.Range('AG2').Value = 2022
$ws.Range('AH2').Value = "'"
$ws.Range('AI2').Value = "'"
$ws.Range('AJ2').Formula = '=HYPERLINK("CertificateUpload_1688981249563.pdf")'
$ws.Range('AK2').Value = "'kjh"
$ws.Range('AL2').Value = "'kjh"
$ws.Range('AM2').Formula = '=HYPERLINK("InternUpload_1688981249627.pdf")'
$ws.Range('AN2').Value = "'"
$ws.Range('AO2').Value = "'firstyear"
$ws.Range('AP2').Value = "'2019"
$ws.Range('AQ2').Value = "'jaya"
$ws.Range('AR2').Value = "'madhuri,89129"
$ws.Range('AS2').Value = "'andhrauniversity"
$ws.Range('AT2').Value = "'nucluesteq raipur"
$ws.Range('AU2').Value = "'supriya"
$ws.Range('AV2').Value = "'76"
$ws.Range('AW2').Value = "'jhb"
$ws.Range('AX2').Value = "'j"
$ws.Range('AY2').Value = "'"
$ws.Range('AZ2').Value = "'"
$ws.Range('BA2').Formula = '=HYPERLINK("Upload_1689317571425.pdf")'
$ws.Range('BB2').Value = "'"
$ws.Range('BC2').Value = "'"